$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 133; this shifts the existing rows 133-177
# down to 134-178 (matches the diff: old row133 data now lives in row134,
# ..., old row177 data now lives in the newly added row178).
$ws.Rows("133:133").Insert()

# Populate the newly inserted row 133 with the new record's data.
$ws.Range("A133").Value = 11
$ws.Range("B133").Value = "Vega Monumental Concepción"
$ws.Range("C133").Value = "Bíobío"
$ws.Range("D133").Value = 45007
$ws.Range("E133").Value = 8
$ws.Range("F133").Value = 100112021
$ws.Range("G133").Value = "Ají"
$ws.Range("H133").Value = "Americana (o)"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 100
$ws.Range("K133").Value = 20000
$ws.Range("L133").Value = 22000
$ws.Range("M133").Value = 21000
$ws.Range("N133").Value = "`$/saco 25 kilos"
$ws.Range("O133").Value = "Región Metropolitana"
$ws.Range("P133").Value = 840
$ws.Range("Q133").Value = 25
$ws.Range("R133").Value = "Hortaliza"
